$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @{
  2  = "1/05/23"
  3  = "2/05/23"
  4  = "3/05/23"
  5  = "4/05/23"
  6  = "5/05/23"
  7  = "8/05/23"
  8  = "9/05/23"
  9  = "10/05/23"
  10 = "11/05/23"
  11 = "12/05/23"
  12 = "15/05/23"
  13 = "16/05/23"
  14 = "17/05/23"
  15 = "18/05/23"
  16 = "19/05/23"
  17 = "22/05/23"
  18 = "23/05/23"
  19 = "24/05/23"
  20 = "25/05/23"
  21 = "26/05/23"
  22 = "30/05/23"
  23 = "31/05/23"
  24 = "1/05/24"
  25 = "2/05/24"
  26 = "3/05/24"
  27 = "6/05/24"
  28 = "7/05/24"
  29 = "8/05/24"
  30 = "9/05/24"
  31 = "10/05/24"
  32 = "13/05/24"
  33 = "14/05/24"
  34 = "15/05/24"
  35 = "16/05/24"
  36 = "17/05/24"
  37 = "20/05/24"
  38 = "21/05/24"
  39 = "22/05/24"
  40 = "23/05/24"
  41 = "24/05/24"
  42 = "28/05/24"
  43 = "29/05/24"
  44 = "30/05/24"
  45 = "31/05/24"
  46 = "1/05/25"
  47 = "2/05/25"
  48 = "5/05/25"
  49 = "6/05/25"
  50 = "7/05/25"
  51 = "8/05/25"
  52 = "9/05/25"
}

foreach ($row in $dates.Keys) {
  $cell = $ws.Cells.Item($row, 1)
  $cell.NumberFormat = "@"
  $cell.Value = $dates[$row]
  $cell.Style = "Normal"
}
